# Generate Report for Handback
# Row 7 on both the "zh-cn" and "de-de" sheets gets a new handback result:
#   - Latest Target File (I)      : the handback markdown file, now a hyperlink
#   - Latest Handback File (J)    : the generated xlf file name
#   - Latest Handback DateTime (K): the timestamp of this handback attempt
#   - Error Detail (P)            : "not the latest version" error message

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8514ca4b895b09d7b77cbcfdcb7419a0fd6e224/e2e/2983f36a-7166-4fae-a14e-c5478578f842.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7d0fd8918f8009b5c47b3d55b8ba9efa8cd7668/e2e/2983f36a-7166-4fae-a14e-c5478578f842.md."
$handbackMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7d0fd8918f8009b5c47b3d55b8ba9efa8cd7668/e2e/2983f36a-7166-4fae-a14e-c5478578f842.md"
$handbackMdName = "2983f36a-7166-4fae-a14e-c5478578f842.md"

# ---- zh-cn sheet, row 7 ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $handbackMdUrl, "", "", $handbackMdName)
$wsZh.Range("I7").Font.Underline = 2
$wsZh.Range("I7").Font.Color = 15570276
$wsZh.Range("J7").Value = "2983f36a-7166-4fae-a14e-c5478578f842.73da93e233d0b2374e5a70aeccaee15d8389f725.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-03 22:58:55"
$wsZh.Range("P7").Value = $errorMessage

# ---- de-de sheet, row 7 ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $handbackMdUrl, "", "", $handbackMdName)
$wsDe.Range("I7").Font.Underline = 2
$wsDe.Range("I7").Font.Color = 15570276
$wsDe.Range("J7").Value = "2983f36a-7166-4fae-a14e-c5478578f842.73da93e233d0b2374e5a70aeccaee15d8389f725.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-03 22:59:03"
$wsDe.Range("P7").Value = $errorMessage
